$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Save" header + data, matching the existing header style (s=1)
# used by B1:G1. Copy format only from G1 so the same cellXf/style is reused
# rather than a brand-new style being created.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
